# Updated Master data as per 16th May Refresh
# Adds three new rows (34-36) to the reg_center_user_h master data sheet,
# following the same pattern as the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 34; A = 10005; B = 110033 },
    @{ Row = 35; A = 10005; B = 110034 },
    @{ Row = 36; A = 10005; B = 110035 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = "eng"
    $ws.Cells.Item($rowNum, 4).Value = $true
    $ws.Cells.Item($rowNum, 5).Value = "superadmin"
    $ws.Cells.Item($rowNum, 6).Value = "now()"
    $ws.Cells.Item($rowNum, 7).Value = "now()"
}

# Mirror the selection state captured in the saved workbook: selecting the
# empty rows below the new data (as if the user clicked row 37's header).
$ws.Range("A37:XFD1048576").Select()
